$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data entered on the "Requerimiento" detail row (row 15) and the
#     document reference field (G2) ---
$ws.Range("A15").Value = "RF_3.1.1a"
$ws.Range("D15").Value = "Funcional"
$ws.Range("C15").Value = "Baja"
$ws.Range("B15").Value = "Inicial"
$ws.Range("E15").Value = "RF_3.1.1b"
$ws.Range("G2").Value = "SRS-SISQSF"

# --- Formatting: A15 becomes a boxed, justified "ID" style cell; E15 (and
#     the cell below it, A16) switch to the plain Times New Roman style
#     used for free-form entries, losing their thin grid border ---
$a15 = $ws.Range("A15")
$a15.Borders.LineStyle = 1
$a15.Borders.Weight = -4138
$a15.HorizontalAlignment = -4130
$a15.VerticalAlignment = -4108
$a15.WrapText = $true
$a15.Font.Size = 12
$a15.Font.Name = "Times New Roman"
$a15.Font.Family = 1

$e15 = $ws.Range("E15")
$e15.Borders.LineStyle = -4142
$e15.Font.Size = 12
$e15.Font.Name = "Times New Roman"
$e15.Font.Family = 1

$a16 = $ws.Range("A16")
$a16.Borders.LineStyle = -4142
$a16.Font.Size = 12
$a16.Font.Name = "Times New Roman"
$a16.Font.Family = 1

# --- View state: scrolled down a bit with A16 as the active selection ---
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("A16").Select() | Out-Null
